$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) summary sheet: insert a new top row for
#    the freshly added "2022-Q4" quarter, pushing existing rows down.
# ---------------------------------------------------------------------
$ws1.Rows.Item(2).Insert()
$ws1.Range("B2:D2").ClearFormats()
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 0.72
# Give the new index cell (A2) the same style as the other index cells.
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q4" worksheet (same layout as the other
#    quarterly fund-holding sheets) positioned right after "总计".
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q2")
$src.Copy($null, $ws1)
$newWs = $wb.Worksheets.Item("2022-Q2 (2)")
$newWs.Name = "2022-Q4"

# The template sheet only has 5 data rows (rows 2-6); we need 8 data
# rows (rows 2-9), so insert 3 more rows, copying the last row's format.
$newWs.Range("A7:H9").EntireRow.Insert()
$newWs.Range("A6:H6").Copy()
$newWs.Range("A7:H9").PasteSpecial(-4122)

# Columns B:G store fund codes/names/figures as text (not numbers), so
# force text formatting before assigning the numeric-looking strings.
$newWs.Range("B2:G9").NumberFormat = "@"

$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "217024"
$newWs.Range("C2").Value = "招商安盈债券A"
$newWs.Range("D2").Value = "40.95"
$newWs.Range("E2").Value = "20.32"
$newWs.Range("F2").Value = "0.98"
$newWs.Range("G2").Value = "0.4013"
$newWs.Range("H2").Value = 7

$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "014887"
$newWs.Range("C3").Value = "招商安福1年定期开放债券"
$newWs.Range("D3").Value = "17.78"
$newWs.Range("E3").Value = "33.59"
$newWs.Range("F3").Value = "1.11"
$newWs.Range("G3").Value = "0.1974"
$newWs.Range("H3").Value = 10

$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "010430"
$newWs.Range("C4").Value = "招商安阳债券A"
$newWs.Range("D4").Value = "16.90"
$newWs.Range("E4").Value = "20.35"
$newWs.Range("F4").Value = "0.72"
$newWs.Range("G4").Value = "0.1217"
$newWs.Range("H4").Value = 10

$newWs.Range("A5").Value = 3
$newWs.Range("B5").Value = "010431"
$newWs.Range("C5").Value = "招商安阳债券C"
$newWs.Range("D5").Value = "0.11"
$newWs.Range("E5").Value = "20.35"
$newWs.Range("F5").Value = "0.72"
$newWs.Range("G5").Value = "0.0008"
$newWs.Range("H5").Value = 10

$newWs.Range("A6").Value = 4
$newWs.Range("B6").Value = "006857"
$newWs.Range("C6").Value = "蜂巢卓睿灵活配置混合A"
$newWs.Range("D6").Value = "0.07"
$newWs.Range("E6").Value = "68.96"
$newWs.Range("F6").Value = "1.16"
$newWs.Range("G6").Value = "0.0008"
$newWs.Range("H6").Value = 6

$newWs.Range("A7").Value = 5
$newWs.Range("B7").Value = "003366"
$newWs.Range("C7").Value = "浙商汇金中证转型成长指数"
$newWs.Range("D7").Value = "0.06"
$newWs.Range("E7").Value = "93.82"
$newWs.Range("F7").Value = "1.26"
$newWs.Range("G7").Value = "0.0008"
$newWs.Range("H7").Value = 6

$newWs.Range("A8").Value = 6
$newWs.Range("B8").Value = "006858"
$newWs.Range("C8").Value = "蜂巢卓睿灵活配置混合C"
$newWs.Range("D8").Value = "0.03"
$newWs.Range("E8").Value = "68.96"
$newWs.Range("F8").Value = "1.16"
$newWs.Range("G8").Value = "0.0003"
$newWs.Range("H8").Value = 6

$newWs.Range("A9").Value = 7
$newWs.Range("B9").Value = "012233"
$newWs.Range("C9").Value = "招商安盈债券C"
$newWs.Range("D9").Value = "0.01"
$newWs.Range("E9").Value = "20.32"
$newWs.Range("F9").Value = "0.98"
$newWs.Range("G9").Value = "0.0001"
$newWs.Range("H9").Value = 7

# Drop the temporary "@" text number format now that the values are
# safely stored as text, so the cells don't keep a stray style index.
$newWs.Range("B2:G9").NumberFormat = "General"

# ---------------------------------------------------------------------
# 3) Restore the originally active sheet ("总计") as the active tab.
# ---------------------------------------------------------------------
$ws1.Activate()
